# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计", as a copy of the
#    "2021-Q3" sheet so it inherits the same layout/column styles, then
#    overwrite its contents with the new quarter's numbers.
# 2) Update the "总计" (summary) sheet: insert a row for 2022-Q4 and shift
#    the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item("2021-Q3")

# Helper: write $text into $cell as a genuine text value (no numeric
# coercion) without leaving a stray "quote prefix" cell style behind.
# Building a literal string formula and collapsing it back to a value via
# copy / paste-special(values) yields a plain shared-string cell, matching
# how the original file stores these (inlineStr / s, no style override).
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Step 1: create the new "2022-Q4" sheet as a copy of "2021-Q3" ---
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Header: 基金金额 -> 基金规模
Set-TextValue $q4.Cells.Item(1, 4) "基金规模"

# Row 2: 519222 / 海富通欣益灵活配置混合A / 0.25 / 31.65 / 0.16 / 0.0004 / 8
Set-TextValue $q4.Cells.Item(2, 2) "519222"
Set-TextValue $q4.Cells.Item(2, 3) "海富通欣益灵活配置混合A"
Set-TextValue $q4.Cells.Item(2, 4) "0.25"
Set-TextValue $q4.Cells.Item(2, 5) "31.65"
Set-TextValue $q4.Cells.Item(2, 6) "0.16"
Set-TextValue $q4.Cells.Item(2, 7) "0.0004"
$q4.Cells.Item(2, 8).Value = 8

# Row 3: 519221 / 海富通欣益灵活配置混合C / 0.10 / 31.65 / 0.16 / 0.0002 / 8
Set-TextValue $q4.Cells.Item(3, 2) "519221"
Set-TextValue $q4.Cells.Item(3, 3) "海富通欣益灵活配置混合C"
Set-TextValue $q4.Cells.Item(3, 4) "0.10"
Set-TextValue $q4.Cells.Item(3, 5) "31.65"
Set-TextValue $q4.Cells.Item(3, 6) "0.16"
Set-TextValue $q4.Cells.Item(3, 7) "0.0002"
$q4.Cells.Item(3, 8).Value = 8

# --- Step 2: update "总计" sheet, shifting rows 2-4 down to 3-5 and ---
# --- inserting the new 2022-Q4 row at row 2 ---

# Bring A5 into existence with the same style as A4 (Range.Copy carries
# both value and formatting), then overwrite the values bottom-up so we
# never clobber data we still need to read.
$total.Range("A4").Copy($total.Range("A5"))

$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2021-Q1"
$total.Cells.Item(5, 3).Value = 4
$total.Cells.Item(5, 4).Value = 0.1

$total.Cells.Item(4, 2).Value = "2021-Q2"
$total.Cells.Item(4, 3).Value = 6
$total.Cells.Item(4, 4).Value = 0.16

$total.Cells.Item(3, 2).Value = "2021-Q3"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.04

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0
